$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.356.03'
$ws.Range('E2').Value = '  +5.18%  '
$ws.Range('D3').Value = '3.514.67'
$ws.Range('E3').Value = '  +2.89%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '''419.13'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').Value = '''132.98'
$ws.Range('E6').Value = '  +3.34%  '
$ws.Range('D7').Value = '''0.656'
$ws.Range('E7').Value = '  +5.04%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '''0.781'
$ws.Range('E9').Value = '  +7.52%  '
$ws.Range('E10').Value = '  +17.06%  '
$ws.Range('D11').Value = '''43.70'
$ws.Range('E11').Value = '  +2.29%  '
$ws.Range('D12').Value = '''0.0000268'
$ws.Range('E12').Value = '  +23.88%  '
$ws.Range('D13').Value = '''10.10'
$ws.Range('E13').Value = '  +9.79%  '
$ws.Range('D14').Value = '4.071.84'
$ws.Range('E14').Value = '  +2.66%  '
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D16').Value = '''20.65'
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('D17').Value = '3.530.04'
$ws.Range('E17').Value = '  +3.53%  '
$ws.Range('D18').Value = '''12.82'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').Value = '''1.11'
$ws.Range('E19').Value = '  +3.95%  '
$ws.Range('D20').Value = '65.266.40'
$ws.Range('E20').Value = '  +4.95%  '
$ws.Range('D21').Value = '''456.40'
$ws.Range('E21').Value = '  -4.09%  '
$ws.Range('D22').Value = '''90.60'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('D24').Value = '''13.42'
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('D25').Value = '''3.42'
$ws.Range('E25').Value = '  +4.39%  '
$ws.Range('D26').Value = '''10.01'
$ws.Range('E26').Value = '  +3.03%  '
$ws.Range('D27').Value = '''34.36'
$ws.Range('E27').Value = '  +3.20%  '
$ws.Range('D28').Value = '''12.79'
$ws.Range('E28').Value = '  +7.95%  '
$ws.Range('D29').Value = '''2.75'
$ws.Range('E29').Value = '  +3.80%  '
$ws.Range('D30').Value = '''7.48'
$ws.Range('E30').Value = '  -3.51%  '
$ws.Range('E31').Value = '  +6.05%  '
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('D33').Value = '''40.01'
$ws.Range('E33').Value = '  -2.33%  '
$ws.Range('D34').Value = '''0.999'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = '''57.49'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').Value = '''0.0509'
$ws.Range('E36').Value = '  +4.45%  '
$ws.Range('D37').Value = '0.0₃0743'
$ws.Range('E37').Value = '  +40.56%  '
$ws.Range('E38').Value = '  +10.73%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''3.08'
$ws.Range('E39').Value = '  +1.33%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '''0.998'
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').Value = '''4.56'
$ws.Range('E41').Value = '  +5.80%  '
$ws.Range('D42').Value = '''2.75'
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('D43').Value = '''146.34'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('D45').Value = '''0.313'
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('E46').Value = '  -2.54%  '
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('D48').Value = '''15.95'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('D49').Value = '''0.147'
$ws.Range('E49').Value = '  +5.69%  '
$ws.Range('D50').Value = '''2.57'
$ws.Range('E50').Value = '  +11.34%  '
$ws.Range('D51').Value = '''21.67'
$ws.Range('E51').Value = '  -2.41%  '
